# "Update to work live on the server"
# Fill in the self-evaluation scores that were missing/changed and move the
# active selection to reflect where the author was working (C11:C32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Basic Options section (rows 11-32) ---------------------------------
# E11 (Comments for "Login Screen") was empty, now has a value.
$ws.Cells.Item(11, 5).Value = 9

# C13 (Numbers of Commits - "User Home Screen") 5 -> 5.5
$ws.Cells.Item(13, 3).Value = 5.5

# C14 (Numbers of Commits - "Publish New Ad") 5 -> 2.5
$ws.Cells.Item(14, 3).Value = 2.5

# C15 (Numbers of Commits - "Adding Picture to Publish New Ad") 5 -> 7
$ws.Cells.Item(15, 3).Value = 7

# E18 (Comments for "Implement Paging") was empty, now has a value.
$ws.Cells.Item(18, 5).Value = 4

# --- Admin Options section (rows 34-50) ---------------------------------
# C37 (Numbers of Commits - "Admin Approve Ad") empty -> 0.2
$ws.Cells.Item(37, 3).Value = 0.2

# C38 (Numbers of Commits - "Admin Reject Ad") empty -> 0.2
$ws.Cells.Item(38, 3).Value = 0.2

# The Total Score in C51 is =SUM(C11:C50); it recalculates automatically.

# --- Selection / view state ----------------------------------------------
# Move the live selection to where the scores were just edited.
$ws.Range("C11:C32").Select()
